# Refresh Market Board derived profit columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# across each job sheet, per the scheduled Sheets runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1702.3077
$ws.Range("I28").Value = 217.4
$ws.Range("J28").Value = 6652
$ws.Range("K28").Value = 217.4
$ws.Range("L28").Value = 6652
$ws.Range("M28").Value = 267.6
$ws.Range("N28").Value = -7622

$ws.Range("H113").Value = 3150.7144
$ws.Range("I113").Value = 2515.7144
$ws.Range("J113").Value = 3785.7144
$ws.Range("K113").Value = 2515.7144
$ws.Range("L113").Value = 3785.7144
$ws.Range("M113").Value = 738.2856000000002

$ws.Range("H132").Value = 5496372.5
$ws.Range("I132").Value = 5716207
$ws.Range("J132").Value = 506
$ws.Range("K132").Value = 17148621
$ws.Range("L132").Value = 1518
$ws.Range("M132").Value = -17146091
$ws.Range("N132").Value = -6578

$ws.Range("H137").Value = 1150.1428
$ws.Range("I137").Value = 1005.06665
$ws.Range("J137").Value = 2020.6
$ws.Range("K137").Value = 3015.19995
$ws.Range("L137").Value = 6061.799999999999
$ws.Range("M137").Value = -465.1999500000002
$ws.Range("N137").Value = -11161.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3125
$ws.Range("I2").Value = 1500
$ws.Range("J2").Value = 3666.6667
$ws.Range("K2").Value = 1500
$ws.Range("L2").Value = 3666.6667
$ws.Range("M2").Value = -1387
$ws.Range("N2").Value = -3892.6667

$ws.Range("H32").Value = 19282.828
$ws.Range("I32").Value = 21076.553
$ws.Range("J32").Value = 6726.75
$ws.Range("K32").Value = 21076.553
$ws.Range("L32").Value = 6726.75
$ws.Range("M32").Value = -20789.553

$ws.Range("H45").Value = 1077.5
$ws.Range("I45").Value = 1044.4445
$ws.Range("J45").Value = 1137
$ws.Range("K45").Value = 1044.4445
$ws.Range("L45").Value = 1137
$ws.Range("M45").Value = -667.4445000000001

$ws.Range("H61").Value = 2155.8572
$ws.Range("I61").Value = 1131.3334
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 1131.3334
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -919.3334
$ws.Range("N61").Value = -4424

$ws.Range("H74").Value = 1206.3043
$ws.Range("I74").Value = 1374.6428
$ws.Range("J74").Value = 944.44446
$ws.Range("K74").Value = 1374.6428
$ws.Range("L74").Value = 944.44446
$ws.Range("M74").Value = -500.6428000000001
$ws.Range("N74").Value = -2692.44446

$ws.Range("H77").Value = 1206.3043
$ws.Range("I77").Value = 1374.6428
$ws.Range("J77").Value = 944.44446
$ws.Range("K77").Value = 6873.214
$ws.Range("L77").Value = 4722.2223
$ws.Range("M77").Value = -2505.214
$ws.Range("N77").Value = -13458.2223

$ws.Range("H102").Value = 2526
$ws.Range("I102").Value = 2526
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2526
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -904
$ws.Range("N102").ClearContents()

$ws.Range("H110").Value = 1267.0834
$ws.Range("I110").Value = 765.73334
$ws.Range("J110").Value = 2102.6667
$ws.Range("K110").Value = 765.73334
$ws.Range("L110").Value = 2102.6667
$ws.Range("M110").Value = 1279.26666

$ws.Range("H116").Value = 3125
$ws.Range("I116").Value = 1500
$ws.Range("J116").Value = 3666.6667
$ws.Range("K116").Value = 1500
$ws.Range("L116").Value = 3666.6667
$ws.Range("M116").Value = 794
$ws.Range("N116").Value = -8254.6667

$ws.Range("H122").Value = 901.36365
$ws.Range("I122").Value = 879.44446
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2638.33338
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -188.33338

$ws.Range("H132").Value = 4508.8047
$ws.Range("I132").Value = 4951.8
$ws.Range("J132").Value = 3300.6365
$ws.Range("K132").Value = 14855.4
$ws.Range("L132").Value = 9901.9095
$ws.Range("M132").Value = -12325.4
$ws.Range("N132").Value = -14961.9095

$ws.Range("H136").Value = 2155.8572
$ws.Range("I136").Value = 1131.3334
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 3394.0002
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -844.0001999999999
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3125
$ws.Range("I3").Value = 1500
$ws.Range("J3").Value = 3666.6667
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 3666.6667
$ws.Range("M3").Value = -1386
$ws.Range("N3").Value = -3894.6667

$ws.Range("H86").Value = 1636.6154
$ws.Range("I86").Value = 1495.8572
$ws.Range("J86").Value = 1800.8334
$ws.Range("K86").Value = 1495.8572
$ws.Range("L86").Value = 1800.8334
$ws.Range("M86").Value = -372.8571999999999
$ws.Range("N86").Value = -4046.8334

$ws.Range("H89").Value = 1636.6154
$ws.Range("I89").Value = 1495.8572
$ws.Range("J89").Value = 1800.8334
$ws.Range("K89").Value = 7479.286
$ws.Range("L89").Value = 9004.166999999999
$ws.Range("M89").Value = -1863.286
$ws.Range("N89").Value = -20236.167

$ws.Range("H105").Value = 4370.8887
$ws.Range("I105").Value = 4187.077
$ws.Range("J105").Value = 4848.8
$ws.Range("K105").Value = 4187.077
$ws.Range("L105").Value = 4848.8
$ws.Range("M105").Value = -2440.077
$ws.Range("N105").Value = -8342.799999999999

$ws.Range("H134").Value = 39775.406
$ws.Range("I134").Value = 55191.367
$ws.Range("J134").Value = 3162.5
$ws.Range("K134").Value = 165574.101
$ws.Range("L134").Value = 9487.5
$ws.Range("M134").Value = -163039.101

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10003110
$ws.Range("I31").Value = 2454.2144
$ws.Range("J31").Value = 33337974
$ws.Range("K31").Value = 2454.2144
$ws.Range("L31").Value = 33337974
$ws.Range("M31").Value = -2159.2144
$ws.Range("N31").Value = -33338564

$ws.Range("H34").Value = 10003110
$ws.Range("I34").Value = 2454.2144
$ws.Range("J34").Value = 33337974
$ws.Range("K34").Value = 2454.2144
$ws.Range("L34").Value = 33337974
$ws.Range("M34").Value = -2252.2144
$ws.Range("N34").Value = -33338378

$ws.Range("H58").Value = 1528.4
$ws.Range("I58").Value = 1476
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 1476
$ws.Range("L58").Value = 2000
$ws.Range("M58").Value = -1273
$ws.Range("N58").Value = -2406

$ws.Range("H105").Value = 1786.8572
$ws.Range("I105").Value = 1427.25
$ws.Range("J105").Value = 2266.3333
$ws.Range("K105").Value = 1427.25
$ws.Range("L105").Value = 2266.3333
$ws.Range("M105").Value = 319.75
$ws.Range("N105").Value = -5760.3333

$ws.Range("H132").Value = 1961.3077
$ws.Range("I132").Value = 1236.1818
$ws.Range("J132").Value = 5949.5
$ws.Range("K132").Value = 3708.5454
$ws.Range("L132").Value = 17848.5
$ws.Range("M132").Value = -1178.5454

$ws.Range("H134").Value = 1070.28
$ws.Range("I134").Value = 1010.7083
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 3032.1249
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -497.1248999999998
$ws.Range("N134").Value = -12570

$ws.Range("H136").Value = 1528.4
$ws.Range("I136").Value = 1476
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 4428
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -1878
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1259649.1
$ws.Range("I131").Value = 22446
$ws.Range("J131").Value = 1701507.4
$ws.Range("K131").Value = 67338
$ws.Range("L131").Value = 5104522.199999999
$ws.Range("M131").Value = -62298
$ws.Range("N131").Value = -5114602.199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3240.7727
$ws.Range("I122").Value = 3828.1428
$ws.Range("J122").Value = 2966.6667
$ws.Range("K122").Value = 11484.4284
$ws.Range("L122").Value = 8900.000100000001
$ws.Range("M122").Value = -9034.428400000001

$ws.Range("H132").Value = 58595.5
$ws.Range("I132").Value = 69374.03
$ws.Range("J132").Value = 4702.8335
$ws.Range("K132").Value = 208122.09
$ws.Range("L132").Value = 14108.5005
$ws.Range("M132").Value = -205592.09
$ws.Range("N132").Value = -19168.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2833.1667
$ws.Range("I40").Value = 3249.75
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 3249.75
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -3113.75
$ws.Range("N40").Value = -2272

$ws.Range("H122").Value = 5267.3335
$ws.Range("I122").Value = 9802
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 29406
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -26956
$ws.Range("N122").Value = -13900

$ws.Range("H132").Value = 1868.8334
$ws.Range("I132").Value = 1768.35
$ws.Range("J132").Value = 2069.8
$ws.Range("K132").Value = 5305.049999999999
$ws.Range("L132").Value = 6209.400000000001
$ws.Range("M132").Value = -2775.049999999999
$ws.Range("N132").Value = -11269.4

$ws.Range("H136").Value = 5235.7036
$ws.Range("I136").Value = 6261.263
$ws.Range("J136").Value = 2800
$ws.Range("K136").Value = 18783.789
$ws.Range("L136").Value = 8400
$ws.Range("M136").Value = -16233.789
$ws.Range("N136").Value = -13500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2344.4443
$ws.Range("I122").Value = 2928.75
$ws.Range("J122").Value = 1877
$ws.Range("K122").Value = 8786.25
$ws.Range("L122").Value = 5631
$ws.Range("M122").Value = -6336.25
$ws.Range("N122").Value = -10531

$ws.Range("H132").Value = 2252.743
$ws.Range("I132").Value = 1749.72
$ws.Range("J132").Value = 3510.3
$ws.Range("K132").Value = 5249.16
$ws.Range("L132").Value = 10530.9
$ws.Range("M132").Value = -2719.16

$ws.Range("H136").Value = 5271.5864
$ws.Range("I136").Value = 6628.4546
$ws.Range("J136").Value = 1007.1429
$ws.Range("K136").Value = 19885.3638
$ws.Range("L136").Value = 3021.4287
$ws.Range("M136").Value = -17335.3638
$ws.Range("N136").Value = -8121.4287
